$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G30").Value = 2
$ws.Range("G31").Value = 1
$ws.Range("G32").Value = 0
$ws.Range("G33").Value = 0
$ws.Range("G34").Value = 0

$ws.Range("B43").Value = 0

$ws.Range("H32").Select()
